# Apply resume text edits via Word COM-interop Find/Replace.
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.Execute($find, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replace, 2)
}

Replace-Text "Resume - Alex Wilber" "Resume: Alex Wilber"
Replace-Text "Spark Animation: Animation Designer (Jan 2021 - Present)" "Spark Animation: Animation Designer (Jan. 2021 - heute)"
Replace-Text "Pixel Studio: Animations-Designer (Jun 2018 - Dez 2020)" "Pixel Studio: Animation Designer (Jun. 2018 - Dez. 2020)"
Replace-Text "Flash Animation: Junior Animation Designer (Sep 2016 - Mai 2018)" "Flash Animation: Junior Animation Designer (Sept. 2016 - Mai 2018)"
Replace-Text "Master of Arts in Animation, erwartet Abschluss: Dez 2025" "Master of Arts in Animation, erwarteter Abschluss: Dez. 2025"
